# Add the "Eagle_AC" Animal Companion worksheet.
#
# The new sheet is built by copying the existing "Owl_AC" sheet (same
# layout / styles / column widths as the target) to the end of the
# workbook, renaming it, and then overwriting the handful of cells that
# differ for the Eagle (the Listen/Spot note in row 2, the Low-Light
# Vision* special quality in E6, and the Improved Evasion special in F11).

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Owl_AC")
$source.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$eagle = $wb.Worksheets.Item($wb.Worksheets.Count)
$eagle.Name = "Eagle_AC"

# Order matters here so the new shared-string entries land in the same
# order as the target workbook: Low Light Vision* (E6), then the
# Listen/Spot note (A2), then Improved Evasion (F11).
$eagle.Range("E6").Value = "Low Light Vision*"
$eagle.Range("A2").Value = "Listen +4, Spot +16"
$eagle.Range("F11").Value = "Improved Evasion"

# Match the recorded selection/active cell on the new sheet.
$eagle.Range("E9").Select()
